$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 gets overwritten with what used to be row 8's data (the Oct 10 2020
# RCB vs Chennai Super Kings match), and the remaining data rows (3-8) are
# removed so only the header + this single match row remain.
$ws.Range("A2").Value = " Oct 10 2020"
$ws.Range("B2").Value = " Dubai (DSC)"
$ws.Range("C2").Value = "RCB won by 37 runs"
$ws.Range("D2").Value = "Royal Challengers Bangalore"
$ws.Range("E2").Value = "Chennai Super Kings"
# F2 (player name) already reads "Washington Sundar " in both the old and
# new row 2, so it is left untouched.

# Numeric-looking columns must stay text (format code "@") so values like
# "100.00" keep their trailing zeros instead of becoming the number 100.
# The temporary text format is reset back to Normal afterwards so the cell
# style matches the original (unstyled) workbook.
$ws.Range("G2:K2").NumberFormat = "@"
$ws.Range("G2").Value = "10"
$ws.Range("H2").Value = "10"
$ws.Range("I2").Value = "0"
$ws.Range("J2").Value = "1"
$ws.Range("K2").Value = "100.00"
$ws.Range("G2:K2").Style = "Normal"

# Remove the old rows 3-8 entirely (they duplicated data now folded into row 2).
$ws.Range("A3:K8").EntireRow.Delete()
